$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds a header row (row 1) plus two identical data
# rows (rows 2-3) describing a storage/color selection. Katalon AI appended
# two more duplicate data rows (4 and 5) with the same values, extending
# the used range from A1:C3 to A1:C5.

$ws.Cells.Item(4, 1).Value2 = " "
$ws.Cells.Item(4, 2).Value2 = "\31 52174-case-656"
$ws.Cells.Item(4, 3).Value2 = "Black"

$ws.Cells.Item(5, 1).Value2 = " "
$ws.Cells.Item(5, 2).Value2 = "\31 52174-case-656"
$ws.Cells.Item(5, 3).Value2 = "Black"
